# The "metadata" sheet had its fields in a confusing order that tripped up
# the site's parser: the "discord" row sat in the middle of the list
# (row 8) and there was a stray leftover instruction line further down
# (row 12). Fix it by moving "discord" to the end of the field list and
# dropping the stray instruction row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "metadata"
$ws2 = $wb.Worksheets.Item(2)   # "Mid"

# 1. Remove the "discord" row (row 8). Everything below shifts up one row,
#    so the old row 9 (portrait) becomes row 8, old row 10 (description/bio)
#    becomes row 9, old row 11 (extra message) becomes row 10, and the old
#    stray-instruction row 12 becomes row 11.
$ws1.Rows("8").Delete()

# 2. Overwrite that now-row-11 (previously the stray instruction text) with
#    the "discord" entry instead, so it becomes the last field in the list.
$ws1.Range("A11").Value = "discord"
$ws1.Range("B11").Value = "Yassin#8026"
$ws1.Range("C11").Value = "If you aren't sure for portrait, just send me the picture you want on discord (usephysics#0001)"
$ws1.Range("A11").Font.Bold = $true

# 3. The hyperlink that was attached to the portrait URL cell (old B9) needs
#    to follow it to its new address (B8). Rebuild all the hyperlinks on the
#    sheet so they point at the right (post-delete) cells.
$hl = $ws1.Hyperlinks
$hl.Delete()
$hl.Add($ws1.Range("B4"), "https://twitter.com/DJ_Y4SSIN") | Out-Null
$hl.Add($ws1.Range("B5"), "https://www.twitch.tv/dj_y4ssin") | Out-Null
$hl.Add($ws1.Range("B6"), "https://na.op.gg/summoner/userName=DJ%20Y4ssin") | Out-Null
$hl.Add($ws1.Range("B8"), "https://i.imgur.com/lS5ATxi.jpg") | Out-Null
$hl.Add($ws1.Range("B7"), "https://www.youtube.com/channel/UCuYltBAWI35gExZtx1gT3HA") | Out-Null

# 4. The two pictures anchored below row 8 (the metadata guide + example
#    screenshot) need to move up by one row (15pt) to stay under the same
#    cells now that row 8 is gone.
$shapes = $ws1.Shapes
$shapes.Item(1).Top = 210.0
$shapes.Item(2).Top = 208.3655905511811

# 5. Make "metadata" the active/selected sheet (with row 8 selected), and
#    drop the "Mid" sheet's selected-tab flag.
$ws1.Activate()
$ws1.Range("A8:XFD8").Select() | Out-Null
